$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 9
    5  = 3
    6  = 5
    7  = 2
    8  = 3
    9  = 3
    10 = 1
    11 = 7
    12 = 7
    13 = 3
    14 = 2
    15 = 5
    16 = 5
    17 = 4
    18 = 4
    19 = 4
    20 = 6
    21 = 3
    22 = 3
    23 = 5
    24 = 3
    25 = 1
    26 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
